# edit.ps1
# Applies updated loading_percent values (case with 380 kV) to Sheet1
# of the workbook, covering rows 2-25 and columns B,D,E,F,G,H,I,J,L,N,O
# (columns C, K, M remain 0 and are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{ Row = 2; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(15.72110122154908,10.86645300597572,15.85710245078357,27.62890969647131,26.15632835382909,13.35883116385893,23.20505505144129,11.00179161324682,10.42157839338351,17.95041611132247,20.03406003332005) },
    @{ Row = 3; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(15.39775669499676,10.89836383576051,15.88829017821858,27.64339210720159,25.95545547908063,13.37925474599039,23.34204407962968,11.01075709708392,10.24747691145154,17.95019540338838,20.0256363005688) },
    @{ Row = 4; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(15.19704392683539,10.91909334185946,15.90889820399814,27.65941028952811,25.84042488367938,13.39424042917218,23.4303965671893,11.01693260459831,10.13983451567429,17.95199361388906,20.02537553329542) },
    @{ Row = 5; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(15.11480747676658,10.92782716016354,15.91766391973884,27.66772870446966,25.7956885016043,13.40096188592164,23.46747018381384,11.01961843101268,10.09583141608646,17.95321482505318,20.02650583749865) },
    @{ Row = 6; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(15.10112842567191,10.92929472073184,15.91914170214326,27.66921808690353,25.78839067839761,13.40211509256405,23.47369089597121,11.02007465103203,10.08851785296955,17.95344716755871,20.02676823810838) },
    @{ Row = 7; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(15.19593651748291,10.91920996873526,15.90901493106248,27.65951522545534,25.8398128265319,13.39432858878686,23.43089222168294,11.01696814046587,10.13924156763264,17.95200810315272,20.02538576874617) },
    @{ Row = 8; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(15.61012488392778,10.87722047127291,15.86755383655238,27.63242393730417,26.08537533367435,13.36536551669054,23.25141047499483,11.00474408426884,10.36173155992696,17.94993999688025,20.03013733646303) },
    @{ Row = 9; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.40072587711343,10.80386253285911,15.79777813242569,27.63584259337463,26.63037901144229,13.32798142138264,22.93296082623348,10.9860672766546,10.79004098257985,17.96112994543564,20.0783220029252) },
    @{ Row = 10; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.96290179542486,10.7553991819939,15.7534820975524,27.67277006734026,27.06592369426189,13.31235480294112,22.71923642968561,10.97553702052582,11.09714233129641,17.97851072755635,20.13722734445039) },
    @{ Row = 11; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.21352963674509,10.73452234922263,15.73483088786111,27.69700847570504,27.270869506366,13.30781491846194,22.62636221658476,10.97143181716482,11.23466612689214,17.98837589192602,20.16906988584222) },
    @{ Row = 12; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.30761963696663,10.72678429851267,15.72798272967408,27.70725202236663,27.34938458649865,13.30646472490206,22.59181570547146,10.96997514294618,11.28638721100461,17.99239039224723,20.1818470251911) },
    @{ Row = 13; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.28739328007677,10.72844338377855,15.72944807148627,27.70499861508574,27.33243582533153,13.30673911195799,22.59922825055397,10.97028452083847,11.27526465309582,17.99151344586474,20.17906337332498) },
    @{ Row = 14; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.22128731193709,10.73388238056186,15.73426318926566,27.69782990636519,27.27731124664161,13.3076964456544,22.62350758452149,10.97131001774452,11.23892862869833,17.98870059697299,20.17010669855864) },
    @{ Row = 15; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.1806866988249,10.73723572564878,15.73724051369981,27.69357739458954,27.24366170461169,13.30833087472947,22.63846042580346,10.97195089180871,11.21662413923343,17.98701386391728,20.16471390922524) },
    @{ Row = 16; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.94641267967791,10.75678701142815,15.75473108656628,27.67133535038763,27.05266089568297,13.31270315252546,22.72539329000521,10.97581903635584,11.08810727494226,17.97790518501367,20.13524735987655) },
    @{ Row = 17; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.80132408190906,10.76908016865602,15.76584431580943,27.65959301398279,26.93718179428405,13.31604307229713,22.77983609032044,10.97836702928848,11.00867720747127,17.97281693905819,20.11845864112016) },
    @{ Row = 18; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.71739668763436,10.77626097165461,15.77237752702514,27.65353977506753,26.87140718819145,13.31820592326053,22.81155975363322,10.97989708703695,10.96278855326996,17.9700748417575,20.10927763708429) },
    @{ Row = 19; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.68890112433609,10.77871119715397,15.77461383591557,27.65161071818671,26.8492501795541,13.31897976998206,22.82237127055237,10.98042623883921,10.94721799476713,17.96917819281291,20.10625094719841) },
    @{ Row = 20; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.8168189195425,10.76776015012523,15.76464668927905,27.66077052353129,26.9494084174126,13.31566250745844,22.77399818766056,10.97808911813324,11.01715395819801,17.9733395140287,20.12019666923957) },
    @{ Row = 21; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.24072705055156,10.73228027280916,15.73284305443662,27.6999066707562,27.29347866085748,13.30740524416211,22.61635926983292,10.97100615316399,11.2496113985777,17.98951925631935,20.17271803275302) },
    @{ Row = 22; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.51297709957049,10.71006847225568,15.71330824768128,27.73168892013685,27.52359744826172,13.30415892965091,22.51696289419995,10.96694723157335,11.39944122322631,18.00171717074946,20.21123134348308) },
    @{ Row = 23; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(17.36813706041096,10.72183418989965,15.72362021387983,27.71416031839676,27.40032306807642,13.30569497826824,22.56968133851236,10.96906159135323,11.31967951148962,17.99505933657598,20.19029530199828) },
    @{ Row = 24; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.80981529816204,10.76835657754333,15.7651876874184,27.66023599817685,26.94387883714631,13.31583380482975,22.77663618323279,10.97821455868518,11.01332231127518,17.97310268700832,20.11940943892816) },
    @{ Row = 25; Cols = @("B","D","E","F","G","H","I","J","L","N","O"); Vals = @(16.18974855795266,10.82275061843148,15.81542625837454,27.62886571746366,26.47650109825117,13.33601546595847,23.01554179400678,10.99055674280859,10.67531323344015,17.95648408536181,20.06114194036396) }
)

foreach ($rowEntry in $data) {
    $r = $rowEntry.Row
    $cols = $rowEntry.Cols
    $vals = $rowEntry.Vals
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
